# adding averages and more checks
# - Training Dashboard: refresh "LAST UPDATE" (I3:I21) to 16-Sep-2025 and
#   recompute "PERIOD TO EXPIRE" (H3:H21), which drops by 8 days to match.
# - Exam Dashboard: shrink the COMMENTS column and rewrite every remark as
#   "date is valid"; re-point the header font so the dashboards match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------
# Training Dashboard: PERIOD TO EXPIRE (H) / LAST UPDATE (I), rows 3-21
# ---------------------------------------------------------------------
$periodToExpire = @{
    3  = 518;  4  = 588;  5  = 405;  6  = 336;  7  = 503
    8  = 489;  9  = 491;  10 = 539;  11 = 349;  12 = 260
    13 = 489;  14 = -23;  15 = -106; 16 = -34;  17 = -34
    18 = 155;  19 = 313;  20 = 313;  21 = 348
}

foreach ($row in 3..21) {
    $ws1.Range("H$row").Value = $periodToExpire[$row]
    $ws1.Range("I$row").Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------
# Exam Dashboard: narrower COMMENTS column + updated remarks, rows 3-10
# ---------------------------------------------------------------------
$ws2.Columns.Item(5).ColumnWidth = 14.14

foreach ($row in 3..10) {
    $ws2.Range("E$row").Value = "date is valid"
}

# ---------------------------------------------------------------------
# Header styling: bold white text on the dark-blue header fill, and drop
# the oversized 14pt title font so both dashboards share the same bold
# font going forward.
# ---------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $titleFont = $ws.Range("A1").Font
    $titleFont.Size = 11
    $titleFont.Color = 16777215

    $headerRange = $ws.Range("A2", $ws.Cells.Item(2, $ws.UsedRange.Columns.Count))
    $headerFont = $headerRange.Font
    $headerFont.Bold = $true
    $headerFont.Color = 16777215
}
